$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.006.50"
$ws.Range("E2").Value = "  +6.16%  "
$ws.Range("D3").Value = "3.061.61"
$ws.Range("E3").Value = "  +5.98%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  +5.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +10.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "3.059.02"
$ws.Range("E8").Value = "  +5.51%  "
$ws.Range("E9").Value = "  +7.71%  "
$ws.Range("E10").Value = "  +10.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.08"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.00%  "
$ws.Range("E12").Value = "  +12.27%  "
$ws.Range("E13").Value = "  +9.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.17%  "
$ws.Range("D15").Value = "3.559.00"
$ws.Range("E15").Value = "  +5.66%  "
$ws.Range("D16").Value = "64.042.53"
$ws.Range("E16").Value = "  +6.41%  "
$ws.Range("E17").Value = "  +4.22%  "
$ws.Range("D18").Value = "3.059.99"
$ws.Range("E18").Value = "  +5.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "477.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +9.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.679"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.69%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +11.90%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +19.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.29%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +10.74%  "
$ws.Range("E29").Value = "  +7.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.25%  "
$ws.Range("E32").Value = "  +4.57%  "
$ws.Range("E33").Value = "  +10.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.62"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0410"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "446.49"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0811"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +26.79%  "
$ws.Range("D41").Value = "2.967.11"
$ws.Range("E41").Value = "  +4.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.20%  "
$ws.Range("E43").Value = "  +2.38%  "
$ws.Range("E44").Value = "  +9.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.261"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +11.10%  "
$ws.Range("E46").Value = "  +15.00%  "
$ws.Range("E48").Value = "  +8.01%  "
$ws.Range("E49").Value = "  +11.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "117.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.39%  "
$ws.Range("E51").Value = "  +9.78%  "
